$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.81
$ws.Range("G2").Value = 2.22
$ws.Range("H2").Value = 1.83
$ws.Range("I2").Value = 11.5
$ws.Range("J2").Value = 2.8
$ws.Range("K2").Value = 980
$ws.Range("P2").Value = 1.33
$ws.Range("F3").Value = 2.26
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 3.2
$ws.Range("I3").Value = 4.6
$ws.Range("J3").Value = 2.66
$ws.Range("K3").Value = 3.35
$ws.Range("P3").Value = 1.51
$ws.Range("Q3").Value = 2.42
$ws.Range("F4").Value = 1.61
$ws.Range("G4").Value = 1.93
$ws.Range("H4").Value = 4.6
$ws.Range("I4").Value = 11
$ws.Range("J4").Value = 3.5
$ws.Range("K4").Value = 6.2
$ws.Range("P4").Value = 1.75
$ws.Range("Q4").Value = 1.82
$ws.Range("S5").Value = 2
$ws.Range("F6").Value = 2.3
$ws.Range("G6").Value = 3.05
$ws.Range("I6").Value = 3.45
$ws.Range("J6").Value = 3.3
$ws.Range("K6").Value = 5.5
$ws.Range("P6").Value = 2.38
$ws.Range("Q6").Value = 1.48
$ws.Range("H7").Value = 2
$ws.Range("I7").Value = 2.04
$ws.Range("P7").Value = 1.77
$ws.Range("Q7").Value = 2.2
$ws.Range("S7").Value = 4.1
$ws.Range("T7").Value = 1.97
$ws.Range("X7").Value = 13
$ws.Range("Y7").Value = 8.6
$ws.Range("AA7").Value = 34
$ws.Range("AG7").Value = 23
$ws.Range("AM7").Value = 150
$ws.Range("AO7").Value = 22
$ws.Range("G8").Value = 1.21
$ws.Range("P8").Value = 2.6
$ws.Range("Q8").Value = 1.51
$ws.Range("F9").Value = 2.94
$ws.Range("G9").Value = 3
$ws.Range("H9").Value = 2.48
$ws.Range("I9").Value = 2.5
$ws.Range("N9").Value = 4.3
$ws.Range("O9").Value = 1.26
$ws.Range("Y9").Value = 12.5
$ws.Range("AB9").Value = 14
$ws.Range("AF9").Value = 22
$ws.Range("AG9").Value = 13.5
$ws.Range("AK9").Value = 32
$ws.Range("AL9").Value = 40
$ws.Range("AN9").Value = 25
$ws.Range("G10").Value = 3.8
$ws.Range("P10").Value = 2.14
$ws.Range("Y10").Value = 11
$ws.Range("AA10").Value = 26
$ws.Range("AC10").Value = 9
$ws.Range("AJ10").Value = 1000
$ws.Range("AK10").Value = 46
$ws.Range("AL10").Value = 55
$ws.Range("AM10").Value = 1000
$ws.Range("I11").Value = 2.36
$ws.Range("J11").Value = 3.75
$ws.Range("P11").Value = 2.22
$ws.Range("Q11").Value = 1.71
$ws.Range("U11").Value = 2.38
$ws.Range("AB11").Value = 16.5
$ws.Range("AO11").Value = 14.5
$ws.Range("F12").Value = 1.55
$ws.Range("G12").Value = 1.59
$ws.Range("H12").Value = 6.4
$ws.Range("I12").Value = 7
$ws.Range("S12").Value = 2.94
$ws.Range("X12").Value = 20
$ws.Range("AH12").Value = 34
$ws.Range("F13").Value = 2.64
$ws.Range("G13").Value = 2.72
$ws.Range("H13").Value = 2.72
$ws.Range("I13").Value = 2.8
$ws.Range("J13").Value = 3.65
$ws.Range("T13").Value = 1.66
$ws.Range("X13").Value = 19
$ws.Range("Y13").Value = 14
$ws.Range("AD13").Value = 13
$ws.Range("AI13").Value = 980
$ws.Range("AJ13").Value = 980
$ws.Range("AK13").Value = 40
$ws.Range("AL13").Value = 80
$ws.Range("AM13").Value = 75
$ws.Range("X14").Value = 23
$ws.Range("AM14").Value = 80
$ws.Range("AN14").Value = 48
$ws.Range("H15").Value = 7.2
$ws.Range("I15").Value = 8
$ws.Range("K15").Value = 5.7
$ws.Range("P15").Value = 2.82
$ws.Range("Z15").Value = 75
$ws.Range("AA15").Value = 230
$ws.Range("AC15").Value = 14
$ws.Range("AD15").Value = 30
$ws.Range("AE15").Value = 95
$ws.Range("AG15").Value = 11.5
$ws.Range("AO15").Value = 80
$ws.Range("G16").Value = 5.2
$ws.Range("I16").Value = 2.78
$ws.Range("J16").Value = 2.48
$ws.Range("H17").Value = 5.8
$ws.Range("K17").Value = 4.4
$ws.Range("O17").Value = 1.31
$ws.Range("P17").Value = 1.96
$ws.Range("U17").Value = 1.97
$ws.Range("AB17").Value = 8.800000000000001
$ws.Range("AE17").Value = 95
$ws.Range("H18").Value = 4.1
$ws.Range("J18").Value = 3.8
$ws.Range("F19").Value = 2.5
$ws.Range("G19").Value = 2.52
$ws.Range("J19").Value = 3.55
$ws.Range("K19").Value = 3.65
$ws.Range("O19").Value = 1.3
$ws.Range("S19").Value = 3.3
$ws.Range("T19").Value = 1.73
$ws.Range("AA19").Value = 55
$ws.Range("AI19").Value = 46
$ws.Range("AK19").Value = 27
$ws.Range("AO19").Value = 32
$ws.Range("G20").Value = 2.94
$ws.Range("Q20").Value = 1.89
$ws.Range("R20").Value = 1.39
$ws.Range("S20").Value = 3.15
$ws.Range("AB20").Value = 13.5
$ws.Range("AH20").Value = 19.5
$ws.Range("AK20").Value = 36
$ws.Range("AO20").Value = 26
$ws.Range("F21").Value = 2.38
$ws.Range("G21").Value = 2.54
$ws.Range("H21").Value = 3
$ws.Range("J21").Value = 3.4
$ws.Range("Q21").Value = 1.89
$ws.Range("S21").Value = 3.45
$ws.Range("X21").Value = 15.5
$ws.Range("AC21").Value = 8.4
$ws.Range("AD21").Value = 14
$ws.Range("AE21").Value = 38
$ws.Range("AF21").Value = 16.5
$ws.Range("AH21").Value = 17.5
$ws.Range("AM21").Value = 100
$ws.Range("AO21").Value = 34
$ws.Range("I22").Value = 4.4
$ws.Range("J22").Value = 4.1
$ws.Range("K22").Value = 4.4
$ws.Range("P22").Value = 2.3
$ws.Range("Q22").Value = 1.65
$ws.Range("G23").Value = 2.08
$ws.Range("J23").Value = 3.75
$ws.Range("K23").Value = 3.85
$ws.Range("Z23").Value = 30
$ws.Range("AA23").Value = 75
$ws.Range("AE23").Value = 120
$ws.Range("AI23").Value = 210
$ws.Range("AO23").Value = 44
$ws.Range("F24").Value = 2.82
$ws.Range("G24").Value = 2.88
$ws.Range("I24").Value = 2.8
$ws.Range("N24").Value = 3.8
$ws.Range("Q24").Value = 2
$ws.Range("S24").Value = 3.45
$ws.Range("T24").Value = 1.76
$ws.Range("U24").Value = 2.18
$ws.Range("Y24").Value = 11.5
$ws.Range("Z24").Value = 18
$ws.Range("AA24").Value = 42
$ws.Range("AC24").Value = 8
$ws.Range("AD24").Value = 12.5
$ws.Range("AF24").Value = 19.5
$ws.Range("AI24").Value = 980
$ws.Range("AJ24").Value = 980
$ws.Range("AK24").Value = 34
$ws.Range("AL24").Value = 44
$ws.Range("AM24").Value = 95
$ws.Range("AO24").Value = 25
$ws.Range("F25").Value = 2.48
$ws.Range("I25").Value = 2.92
$ws.Range("R25").Value = 1.58
$ws.Range("S25").Value = 2.46
$ws.Range("T25").Value = 1.54
$ws.Range("AA25").Value = 48
$ws.Range("AK25").Value = 28
$ws.Range("AL25").Value = 32
